# Generate Report for Handoff
# "b.md" has now been handed off for both locales (zh-cn and de-de), so
# its status/handoff file/datetime info needs to be refreshed across the
# Overview sheet and the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row 3 corresponds to "b.md"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady
$wsOverview.Range("D3").Value = "2016-23-11 14:23:07"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 corresponds to "b.md"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$zhCnHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"

$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("D3").Value = $zhCnHandoffFile
$wsZhCn.Range("E3").Value = "2016-03-11 14:23:03"

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = $zhCnHandoffFile
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 corresponds to "b.md"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$deDeHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"

$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("D3").Value = $deDeHandoffFile
$wsDeDe.Range("E3").Value = "2016-03-11 14:23:07"

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$D$3') {
        $hl.TextToDisplay = $deDeHandoffFile
    }
}
